$wb = $excel.ActiveWorkbook

# --- Dairy sheet: add Buttermilk / Lassi ---
$dairy = $wb.Worksheets.Item("Dairy")
$dairy.Range("B8").Value = "Buttermilk"
$dairy.Range("B9").Value = "Lassi"
$dairy.Range("B10").Select()

# --- Poultry and Meats sheet: add code column + Pork/Prawn rows ---
$poultry = $wb.Worksheets.Item("Poultry and Meats")
$poultry.Range("A2").Value = "DP01"
$poultry.Range("B2").Value = "Chicken"
$poultry.Range("A3").Value = "DP02"
$poultry.Range("B3").Value = "Egg"
$poultry.Range("A4").Value = "DP03"
$poultry.Range("B4").Value = "Fish – fresh"
$poultry.Range("A5").Value = "DP04"
$poultry.Range("B5").Value = "Fish – sea"
$poultry.Range("A6").Value = "DP05"
$poultry.Range("B6").Value = "Mutton"
$poultry.Range("A7").Value = "DP06"
$poultry.Range("B7").Value = "Pork"
$poultry.Range("A8").Value = "DP07"
$poultry.Range("B8").Value = "Prawn"
$poultry.Range("B17").Select()

# --- Nuts and Legumes sheet: add nuts & legumes list ---
$nuts = $wb.Worksheets.Item("Nuts and Legumes")
$nuts.Range("B2").Value = "Peanut"
$nuts.Range("B3").Value = "Walnut"
$nuts.Range("B4").Value = "Pistachio"
$nuts.Range("B5").Value = "Peas"
$nuts.Range("B6").Value = "Masoor Dal"
$nuts.Range("B7").Value = "Toor Dal"
$nuts.Range("B8").Value = "Chikpeas"
$nuts.Range("B9").Value = "Kidney Bean"
$nuts.Range("B10").Value = "Gram"
$nuts.Range("B11").Value = "Almond"
$nuts.Range("B12").Value = "Cashew"
$nuts.Range("B13").Value = "Moong Dal"
$nuts.Range("B14").Value = "Arhar Dal"
$nuts.Range("B15").Value = "Chana Dal"
$nuts.Range("A2").Select()

# --- Fruits sheet becomes the active tab (selection stays D18) ---
$fruits = $wb.Worksheets.Item("Fruits")
$fruits.Activate()
$fruits.Range("D18").Select()
